$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the header row (row 1), shifting all
# existing data rows (2-34) down to (4-36).
$ws.Rows("2:3").Insert()

# The row-insert copies the header row's formatting onto the new rows'
# B:S cells (bold, no border) - clear that so they match the plain data
# rows below.
$ws.Range("B2:S3").ClearFormats()

# Give column A of the two new rows the same "serial number" look as the
# rest of column A (bold, bordered, centered).
$aStyleRange = $ws.Range("A2:A3")
$aStyleRange.Font.Bold = $true
$aStyleRange.HorizontalAlignment = -4108
$aStyleRange.VerticalAlignment = -4160
$aStyleRange.Borders.LineStyle = 1

# --- Row 2: 江苏通付盾区块链科技有限公司 ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "江苏通付盾区块链科技有限公司"
$ws.Range("C2").Value = "江苏省南京市浦口区"
$ws.Range("D2").Value = "电力信息"
$ws.Range("E2").Value = "Java"
$ws.Range("F2").Value = "电力信息部是驻场外包 时间不定 加班多 本部的话9-18点"
$ws.Range("G2").Value = "1h"
$ws.Range("H2").Value = "驻场外包加班多 就一个餐补25 可以调休"
$ws.Range("I2").Value = "全额的12%"
$ws.Range("J2").Value = "看公司情况 一般有"
$ws.Range("K2").Value = "三个月, 8折"
$ws.Range("L2").Value = "网吧工位，外包是戴尔笔记本"
$ws.Range("M2").Value = "5天"
$ws.Range("N2").Value = "我在驻场外包是看客户方"
$ws.Range("Q2").Value = "2022-06-23 10:05:06"

# --- Row 3: 鱼快创领（南京） ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "鱼快创领（南京）"
$ws.Range("C3").Value = "南京市九龙湖企业总部园"
$ws.Range("D3").Value = "硬件"
$ws.Range("E3").Value = "Java"
$ws.Range("F3").Value = "09:00-18:00"
$ws.Range("G3").Value = "1.5h"
$ws.Range("H3").Value = "双休，平时加班不多，看部门"
$ws.Range("I3").Value = "全额的12%"
$ws.Range("J3").Value = "1个月，根据绩效浮动"
$ws.Range("K3").Value = "6个月, 不打折"
$ws.Range("L3").Value = "网吧工位，笔记本+显示器"
$ws.Range("M3").Value = "5天，一般会多送些"
$ws.Range("N3").Value = "钉钉打卡"
$ws.Range("O3").Value = " 福利不错，零食饮料都有"
$ws.Range("Q3").Value = "2022-06-23 10:04:35"
